$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values for rows 2-12 (columns B, C, D)
$data = @(
    @(10.02, 6, 12),
    @(7.907, 6, 12),
    @(5.219000000000003, 5, 11),
    @(8.0619999999999994, 7, 12),
    @(8.9600000000000009, 8, 15),
    @(11.336, 5, 12),
    @(5.2039999999999997, 6, 13),
    @(11.102, 7, 14),
    @(6.8949999999999996, 7, 13),
    @(4.5679999999999996, 7, 13),
    @(4.8070000000000004, 8, 13)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
}

# Delete the now-unused rows 13-17 (A=12..16 records removed)
$ws.Range("A13:D17").EntireRow.Delete()

# Update the zoom scale of the sheet view
$ws.Activate()
$excel.ActiveWindow.Zoom = 115
